$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) onto the new
# header cells so AD1:AF1 match the look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the team record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (W/L/T) is constant for every player row on this sheet.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 90
    $ws.Cells.Item($r, 32).Value = 0
}
